{"js": "const pairs = [\n  [\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\"],\n  [\"769\u00d73=\", \"346\u00d78=\"],\n  [\"571\u00d76=\", \"262\u00d73=\"],\n  [\"930\u00d79=\", \"739\u00d78=\"],\n  [\"211\u00d79=\", \"576\u00d75=\"],\n  [\"410\u00d75=\", \"702\u00d79=\"],\n  [\"136\u00d73=\", \"817\u00d77=\"],\n  [\"243\u00d74=\", \"417\u00d77=\"],\n  [\"214\u00d78=\", \"584\u00d72=\"],\n  [\"292\u00d74=\", \"513\u00d75=\"],\n  [\"413\u00d79=\", \"580\u00d74=\"],\n  [\"171\u00d76=\", \"978\u00d77=\"],\n  [\"873\u00d78=\", \"742\u00d72=\"],\n  [\"572\u00d76=\", \"437\u00d75=\"],\n  [\"157\u00d78=\", \"961\u00d72=\"],\n  [\"117\u00d74=\", \"151\u00d77=\"],\n  [\"608\u00d79=\", \"605\u00d75=\"],\n  [\"617\u00d75=\", \"212\u00d74=\"],\n  [\"572\u00d78=\", \"333\u00d73=\"],\n  [\"919\u00d76=\", \"182\u00d74=\"],\n  [\"716\u00d75=\", \"215\u00d79=\"],\n  [\"661\u00d77=\", \"810\u00d79=\"],\n  [\"382\u00d79=\", \"908\u00d76=\"],\n  [\"544\u00d78=\", \"684\u00d72=\"],\n  [\"531\u00d72=\", \"458\u00d79=\"],\n  [\"926\u00d76=\", \"909\u00d73=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\")\n    ,@(\"769\u00d73=\", \"346\u00d78=\")\n    ,@(\"571\u00d76=\", \"262\u00d73=\")\n    ,@(\"930\u00d79=\", \"739\u00d78=\")\n    ,@(\"211\u00d79=\", \"576\u00d75=\")\n    ,@(\"410\u00d75=\", \"702\u00d79=\")\n    ,@(\"136\u00d73=\", \"817\u00d77=\")\n    ,@(\"243\u00d74=\", \"417\u00d77=\")\n    ,@(\"214\u00d78=\", \"584\u00d72=\")\n    ,@(\"292\u00d74=\", \"513\u00d75=\")\n    ,@(\"413\u00d79=\", \"580\u00d74=\")\n    ,@(\"171\u00d76=\", \"978\u00d77=\")\n    ,@(\"873\u00d78=\", \"742\u00d72=\")\n    ,@(\"572\u00d76=\", \"437\u00d75=\")\n    ,@(\"157\u00d78=\", \"961\u00d72=\")\n    ,@(\"117\u00d74=\", \"151\u00d77=\")\n    ,@(\"608\u00d79=\", \"605\u00d75=\")\n    ,@(\"617\u00d75=\", \"212\u00d74=\")\n    ,@(\"572\u00d78=\", \"333\u00d73=\")\n    ,@(\"919\u00d76=\", \"182\u00d74=\")\n    ,@(\"716\u00d75=\", \"215\u00d79=\")\n    ,@(\"661\u00d77=\", \"810\u00d79=\")\n    ,@(\"382\u00d79=\", \"908\u00d76=\")\n    ,@(\"544\u00d78=\", \"684\u00d72=\")\n    ,@(\"531\u00d72=\", \"458\u00d79=\")\n    ,@(\"926\u00d76=\", \"909\u00d73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $ok = $find.Execute(\n        [ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false,\n        [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2\n    )\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
